$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells that become empty (clear existing values)
$clearCells = @("A1","C1","D1","F2","F3","B4","G4","C5","I5","A6","E6","B7","C8","E9")
foreach ($cell in $clearCells) {
    $ws.Range($cell).ClearContents()
}

# Cells that get new values
$ws.Range("F1").Value = 6
$ws.Range("I1").Value = 3
$ws.Range("A2").Value = 2
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 3
$ws.Range("H2").Value = 4
$ws.Range("A3").Value = 1
$ws.Range("C3").Value = 6
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 7
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 2
$ws.Range("H5").Value = 5
$ws.Range("H6").Value = 8
$ws.Range("I6").Value = 1
$ws.Range("A7").Value = 7
$ws.Range("D7").Value = 6
$ws.Range("I7").Value = 2
$ws.Range("B8").Value = 4
$ws.Range("E8").Value = 2
$ws.Range("G8").Value = 8
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 2
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 4
